# Logged Week 17 data and fixed Simulate_Season.py tiebreaking method
#
# - Rushing sheet: N.Harris and B.Snell get updated week totals; a new
#   row for D.Watt (a RB who already had receiving stats) is inserted
#   between A.McFarland and D.Johnson, shifting the rest of the roster
#   down by one row and re-numbering the index column.
# - Receiving sheet: several players' totals are updated for the week;
#   no rows are added/removed on this sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Rushing sheet
# ---------------------------------------------------------------------
$rushing = $wb.Worksheets.Item("Rushing")

# Week 17 totals for already-listed players.
$rushing.Range("C4").Value = 171
$rushing.Range("D4").Value = 101
$rushing.Range("E4").Value = 26
$rushing.Range("F4").Value = 28

$rushing.Range("D5").Value = 8

# Insert a new row for D.Watt right after A.McFarland (row 7) and before
# D.Johnson (old row 8), pushing everything below down by one.
$rushing.Rows.Item(8).Insert()

# Match the look of the other index cells (bold, centered/top, boxed)
# instead of whatever default formatting Insert() applied.
$newIndexCell = $rushing.Range("A8")
$newIndexCell.Font.Bold = $true
$newIndexCell.HorizontalAlignment = -4108
$newIndexCell.VerticalAlignment = -4160
$newIndexCell.Borders.LineStyle = 1

# Re-write row 8 through the new last row (13) explicitly so both the
# re-numbered index column and the shifted player rows land correctly.
$rushing.Range("A8").Value = 6
$rushing.Range("B8").Value = "D.Watt"
$rushing.Range("C8").Value = 0
$rushing.Range("D8").Value = 0
$rushing.Range("E8").Value = 1
$rushing.Range("F8").Value = 0

$rushing.Range("A9").Value = 7
$rushing.Range("B9").Value = "D.Johnson"
$rushing.Range("C9").Value = 2
$rushing.Range("D9").Value = 3
$rushing.Range("E9").Value = 0
$rushing.Range("F9").Value = 0

$rushing.Range("A10").Value = 8
$rushing.Range("B10").Value = "C.Claypool"
$rushing.Range("C10").Value = 6
$rushing.Range("D10").Value = 4
$rushing.Range("E10").Value = 1
$rushing.Range("F10").Value = 2

$rushing.Range("A11").Value = 9
$rushing.Range("B11").Value = "J.Washington"
$rushing.Range("C11").Value = 0
$rushing.Range("D11").Value = 2
$rushing.Range("E11").Value = 0
$rushing.Range("F11").Value = 0

$rushing.Range("A12").Value = 10
$rushing.Range("B12").Value = "R.McCloud"
$rushing.Range("C12").Value = 1
$rushing.Range("D12").Value = 1
$rushing.Range("E12").Value = 0
$rushing.Range("F12").Value = 1

$rushing.Range("A13").Value = 11
$rushing.Range("B13").Value = "E.Ebron"
$rushing.Range("C13").Value = 0
$rushing.Range("D13").Value = 0
$rushing.Range("E13").Value = 1
$rushing.Range("F13").Value = 1

# ---------------------------------------------------------------------
# Receiving sheet
# ---------------------------------------------------------------------
$receiving = $wb.Worksheets.Item("Receiving")

$receiving.Range("C2").Value = 87
$receiving.Range("D2").Value = 67

$receiving.Range("C7").Value = 126
$receiving.Range("D7").Value = 87
$receiving.Range("E7").Value = 38
$receiving.Range("G7").Value = 20
$receiving.Range("H7").Value = 10

$receiving.Range("C8").Value = 70
$receiving.Range("D8").Value = 42
$receiving.Range("E8").Value = 30
$receiving.Range("G8").Value = 11

$receiving.Range("E9").Value = 12

$receiving.Range("C10").Value = 50
$receiving.Range("D10").Value = 33
$receiving.Range("G10").Value = 10
$receiving.Range("H10").Value = 3

$receiving.Range("C14").Value = 64
$receiving.Range("D14").Value = 51
$receiving.Range("G14").Value = 18
$receiving.Range("H14").Value = 12

$receiving.Range("C15").Value = 17
$receiving.Range("D15").Value = 14
$receiving.Range("G15").Value = 3
